# Auto-generated edit script: updates cached market-price / profit
# figures on several Leve-profit sheets (refreshed by the scheduled
# market-data runner). Values only -- no formulas are present in
# the source workbook, so we write literals via Range.Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.2143
$ws.Range("I33").Value = 297.75
$ws.Range("K33").Value = 297.75
$ws.Range("M33").Value = -68.75
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H98").Value = 1871.3478
$ws.Range("I98").Value = 1906.4736
$ws.Range("J98").Value = 1704.5
$ws.Range("K98").Value = 1906.4736
$ws.Range("L98").Value = 1704.5
$ws.Range("M98").Value = -408.4736
$ws.Range("N98").Value = -4700.5
$ws.Range("H106").Value = 3153.6667
$ws.Range("I106").Value = 2993
$ws.Range("K106").Value = 2993
$ws.Range("M106").Value = -2362
$ws.Range("H122").Value = 1871.3478
$ws.Range("I122").Value = 1906.4736
$ws.Range("J122").Value = 1704.5
$ws.Range("K122").Value = 5719.4208
$ws.Range("L122").Value = 5113.5
$ws.Range("M122").Value = -3269.4208
$ws.Range("N122").Value = -10013.5
$ws.Range("H127").Value = 2080.3572
$ws.Range("I127").Value = 1177.0834
$ws.Range("K127").Value = 3531.2502
$ws.Range("M127").Value = 1428.7498
$ws.Range("H129").Value = 1712.3636
$ws.Range("J129").Value = 2756.2
$ws.Range("L129").Value = 8268.599999999999
$ws.Range("N129").Value = -18268.6
$ws.Range("H135").Value = 2702.476
$ws.Range("I135").Value = 2059.5
$ws.Range("J135").Value = 4760
$ws.Range("K135").Value = 18535.5
$ws.Range("L135").Value = 42840
$ws.Range("M135").Value = -16000.5
$ws.Range("N135").Value = -47910
$ws.Range("H138").Value = 2995.3333
$ws.Range("J138").Value = 3233.5625
$ws.Range("L138").Value = 9700.6875
$ws.Range("N138").Value = -19980.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24413.965
$ws.Range("I32").Value = 23080.441
$ws.Range("K32").Value = 23080.441
$ws.Range("M32").Value = -22793.441
$ws.Range("H44").Value = 59849.285
$ws.Range("J44").Value = 59849.285
$ws.Range("L44").Value = 59849.285
$ws.Range("N44").Value = -60825.285
$ws.Range("H55").Value = 24551
$ws.Range("J55").Value = 24551
$ws.Range("L55").Value = 24551
$ws.Range("N55").Value = -25181
$ws.Range("H61").Value = 9084.064
$ws.Range("I61").Value = 7415.85
$ws.Range("J61").Value = 12117.182
$ws.Range("K61").Value = 7415.85
$ws.Range("L61").Value = 12117.182
$ws.Range("M61").Value = -7203.85
$ws.Range("N61").Value = -12541.182
$ws.Range("H63").Value = 6931.39
$ws.Range("I63").Value = 3827.25
$ws.Range("J63").Value = 7266.973
$ws.Range("K63").Value = 3827.25
$ws.Range("L63").Value = 7266.973
$ws.Range("M63").Value = -3141.25
$ws.Range("N63").Value = -8638.973
$ws.Range("H66").Value = 6931.39
$ws.Range("I66").Value = 3827.25
$ws.Range("J66").Value = 7266.973
$ws.Range("K66").Value = 19136.25
$ws.Range("L66").Value = 36334.865
$ws.Range("M66").Value = -15704.25
$ws.Range("N66").Value = -43198.865
$ws.Range("H74").Value = 3127.8948
$ws.Range("I74").Value = 1064.8889
$ws.Range("K74").Value = 1064.8889
$ws.Range("M74").Value = -190.8888999999999
$ws.Range("H77").Value = 3127.8948
$ws.Range("I77").Value = 1064.8889
$ws.Range("K77").Value = 5324.4445
$ws.Range("M77").Value = -956.4444999999996
$ws.Range("H97").Value = 1129292.1
$ws.Range("I97").Value = 1691431.4
$ws.Range("J97").Value = 5013.636
$ws.Range("K97").Value = 1691431.4
$ws.Range("L97").Value = 5013.636
$ws.Range("M97").Value = -1690935.4
$ws.Range("N97").Value = -6005.636
$ws.Range("H132").Value = 6029.1143
$ws.Range("I132").Value = 4217.96
$ws.Range("K132").Value = 12653.88
$ws.Range("M132").Value = -10123.88
$ws.Range("H136").Value = 9084.064
$ws.Range("I136").Value = 7415.85
$ws.Range("J136").Value = 12117.182
$ws.Range("K136").Value = 22247.55
$ws.Range("L136").Value = 36351.546
$ws.Range("M136").Value = -19697.55
$ws.Range("N136").Value = -41451.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 92112.37
$ws.Range("I86").Value = 1191.7693
$ws.Range("J86").Value = 223442.11
$ws.Range("K86").Value = 1191.7693
$ws.Range("L86").Value = 223442.11
$ws.Range("M86").Value = -68.76929999999993
$ws.Range("N86").Value = -225688.11
$ws.Range("H89").Value = 92112.37
$ws.Range("I89").Value = 1191.7693
$ws.Range("J89").Value = 223442.11
$ws.Range("K89").Value = 5958.8465
$ws.Range("L89").Value = 1117210.55
$ws.Range("M89").Value = -342.8464999999997
$ws.Range("N89").Value = -1128442.55
$ws.Range("H134").Value = 3370.5293
$ws.Range("I134").Value = 1739.0526
$ws.Range("J134").Value = 8139.4614
$ws.Range("K134").Value = 5217.1578
$ws.Range("L134").Value = 24418.3842
$ws.Range("M134").Value = -2682.1578
$ws.Range("N134").Value = -29488.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2028.6
$ws.Range("I105").Value = 1305
$ws.Range("J105").Value = 2511
$ws.Range("K105").Value = 1305
$ws.Range("L105").Value = 2511
$ws.Range("M105").Value = 442
$ws.Range("N105").Value = -6005
$ws.Range("H132").Value = 27932.477
$ws.Range("I132").Value = 784.8889
$ws.Range("K132").Value = 2354.6667
$ws.Range("M132").Value = 175.3332999999998
$ws.Range("H134").Value = 5357.364
$ws.Range("I134").Value = 4579.222
$ws.Range("K134").Value = 13737.666
$ws.Range("M134").Value = -11202.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1499
$ws.Range("I68").Value = 1665.3334
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 4996.0002
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -4185.0002
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 1499
$ws.Range("I71").Value = 1665.3334
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 14988.0006
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -10932.0006
$ws.Range("N71").Value = -17112
$ws.Range("H140").Value = 1045.4117
$ws.Range("I140").Value = 752.63635
$ws.Range("J140").Value = 1582.1666
$ws.Range("K140").Value = 2257.90905
$ws.Range("L140").Value = 4746.4998
$ws.Range("M140").Value = 2922.09095
$ws.Range("N140").Value = -15106.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 25874.5
$ws.Range("I40").Value = 15000
$ws.Range("K40").Value = 15000
$ws.Range("M40").Value = -14849
$ws.Range("H126").Value = 3388.842
$ws.Range("I126").Value = 2379.5
$ws.Range("J126").Value = 4510.3335
$ws.Range("K126").Value = 7138.5
$ws.Range("L126").Value = 13531.0005
$ws.Range("M126").Value = -4668.5
$ws.Range("N126").Value = -18471.0005
$ws.Range("H132").Value = 4661.5
$ws.Range("I132").Value = 2441.1
$ws.Range("K132").Value = 7323.299999999999
$ws.Range("M132").Value = -4793.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14283.64
$ws.Range("I40").Value = 12673.066
$ws.Range("K40").Value = 12673.066
$ws.Range("M40").Value = -12537.066
$ws.Range("H46").Value = 8515.857
$ws.Range("J46").Value = 9685.111000000001
$ws.Range("L46").Value = 9685.111000000001
$ws.Range("N46").Value = -10061.111
$ws.Range("H68").Value = 3699.889
$ws.Range("I68").Value = 3874.75
$ws.Range("J68").Value = 3560
$ws.Range("K68").Value = 3874.75
$ws.Range("L68").Value = 3560
$ws.Range("M68").Value = -3125.75
$ws.Range("N68").Value = -5058
$ws.Range("H71").Value = 3699.889
$ws.Range("I71").Value = 3874.75
$ws.Range("J71").Value = 3560
$ws.Range("K71").Value = 19373.75
$ws.Range("L71").Value = 17800
$ws.Range("M71").Value = -15629.75
$ws.Range("N71").Value = -25288
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H132").Value = 6688
$ws.Range("I132").Value = 4631
$ws.Range("K132").Value = 13893
$ws.Range("M132").Value = -11363
$ws.Range("H136").Value = 4816.1064
$ws.Range("I136").Value = 4185.278
$ws.Range("J136").Value = 6880.636
$ws.Range("K136").Value = 12555.834
$ws.Range("L136").Value = 20641.908
$ws.Range("M136").Value = -10005.834
$ws.Range("N136").Value = -25741.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4459.6294
$ws.Range("J81").Value = 6699.5
$ws.Range("L81").Value = 13399
$ws.Range("N81").Value = -15521
$ws.Range("H84").Value = 4459.6294
$ws.Range("J84").Value = 6699.5
$ws.Range("L84").Value = 66995
$ws.Range("N84").Value = -77603
$ws.Range("H122").Value = 4999.857
$ws.Range("I122").Value = 4999.857
$ws.Range("K122").Value = 14999.571
$ws.Range("M122").Value = -12549.571
$ws.Range("H136").Value = 3899.05
$ws.Range("I136").Value = 1570.0714
$ws.Range("J136").Value = 9333.333000000001
$ws.Range("K136").Value = 4710.2142
$ws.Range("L136").Value = 27999.999
$ws.Range("M136").Value = -2160.2142
$ws.Range("N136").Value = -33099.999
